$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Window Number"
$ws.Range("B1").Value = "Start (min since recording start time)"
$ws.Range("C1").Value = "Elapsed Time (min)"
$ws.Range("D1").Value = "# Blocks to Divide into"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 5

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 2

# Column widths (closest achievable values given COM ColumnWidth rounding)
$ws.Columns.Item(1).ColumnWidth = 15.16666666666667
$ws.Columns.Item(2).ColumnWidth = 32.66666666666667
$ws.Columns.Item(3).ColumnWidth = 17.33333333333333
$ws.Columns.Item(4).ColumnWidth = 19.33333333333333

# Selection moves to A4 (as if user pressed Enter after filling data)
$ws.Range("A4").Select()
